$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "291.03") need NumberFormat
# forced to Text first, otherwise Excel auto-converts the assigned string
# into a real number -- the source data keeps these as literal text cells.
# NumberFormat is restored immediately after via Style="Normal" so no stray
# cell-style is left behind (matches original General-formatted cells).
$forceTextCells = @("D5", "D6", "D10", "D12", "D14", "D16", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D48")
foreach ($c in $forceTextCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "39.720.11"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.216.83"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "291.03"
$ws.Range("D6").Value = "86.77"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "30.46"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "49.87"
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "2.560.14"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "13.76"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "2.229.55"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "39.686.22"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "5.74"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "65.61"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "237.01"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "23.02"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "9.22"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "156.65"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").Value = "2.97"
$ws.Range("E35").Value = "  +7.15%  "
$ws.Range("D36").Value = "0.0712"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "0.0989"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").Value = "1.74"
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("D41").Value = "15.32"
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").Value = "2.108.99"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "18.16"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "0.0269"
$ws.Range("D46").Value = "9.89"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "1.97"
$ws.Range("E47").Value = "  -7.83%  "
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").Value = "2.434.87"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("E51").Value = "  +2.38%  "

foreach ($c in $forceTextCells) {
    $ws.Range($c).Style = "Normal"
}
